$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values per repull of data / mean calculation fix
$ws.Range("F2").Value = -15
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -6
$ws.Range("F11").Value = -2
